# Apply updated dSF (column F) values to reflect repulled/recalculated data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    8  = 5
    11 = -2
    12 = -3
    18 = -10
    20 = 0
    22 = 5
    25 = 6
    27 = -7
    28 = -1
    30 = 5
    32 = -3
    34 = 9
    35 = -8
    39 = -3
    41 = 3
    42 = -3
    43 = -5
    46 = 5
    50 = -5
    52 = -4
    56 = -6
    60 = -1
    62 = -1
    66 = -7
}

foreach ($row in $updates.Keys) {
    $ws.Cells.Item($row, 6).Value = $updates[$row]
}
